$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of run data appended to the bottom of the "V4" tracking sheet
# (8-8 finished, some of 8-F2 done).
$rows = @(
    @{ Row = 192; A = "Checkpoint 1633";        B = 59432; C = 69850 },
    @{ Row = 193; A = "Checkpoint 1946/1944";    B = 59610; C = 70031 },
    @{ Row = 194; A = "Checkpoint 2388/2385";    B = 59716; C = 70137 },
    @{ Row = 195; A = "Checkpoint 2896/2895";    B = 59863; C = 70285 },
    @{ Row = 196; A = "Get flag";                B = 60022; C = 70444 },
    @{ Row = 197; A = "End Level";               B = 60540; C = 70962 },
    @{ Row = 198; A = "Enter 8-F2";              B = 60909; C = 71707 },
    @{ Row = 199; A = "1st Move";                B = 61137; C = 71957 },
    @{ Row = 200; A = "Platform 1st Move";       B = 61238; C = 72069 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Formula = "=IF(B$row>0,C$row-B$row,0)"
}

# Move the active selection to right after the new last row, matching
# the author's saved cursor position.
[void]$ws.Range("B201").Select()
